# Emergency fix: correct typos / placeholder data in the planning template.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "cliente" (K2) placeholder "HOLA" -> real client name "Brinks"
$ws.Range("K2").Value = "Brinks"

# "ticket_id" (B2) was mistakenly 123 -> correct value is 12
$ws.Range("B2").Value = 12

# Restore the real working selection (J7) instead of the stray I10 selection
$ws.Range("J7").Select()
